$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 245.1875
$ws.Range("I33").Value = 275.6154
$ws.Range("K33").Value = 275.6154
$ws.Range("M33").Value = -46.61540000000002
$ws.Range("H42").Value = 2130.125
$ws.Range("I42").Value = 47
$ws.Range("J42").Value = 2824.5
$ws.Range("K42").Value = 141
$ws.Range("L42").Value = 8473.5
$ws.Range("M42").Value = 89
$ws.Range("N42").Value = -8933.5
$ws.Range("H70").Value = 1744658.8
$ws.Range("I70").Value = 2441324.2
$ws.Range("K70").Value = 7323972.600000001
$ws.Range("M70").Value = -7323702.600000001
$ws.Range("H73").Value = 1744658.8
$ws.Range("I73").Value = 2441324.2
$ws.Range("K73").Value = 7323972.600000001
$ws.Range("M73").Value = -7323036.600000001
$ws.Range("H82").Value = 9382.272000000001
$ws.Range("I82").Value = 7646.3335
$ws.Range("K82").Value = 22939.0005
$ws.Range("M82").Value = -22533.0005
$ws.Range("H85").Value = 9382.272000000001
$ws.Range("I85").Value = 7646.3335
$ws.Range("K85").Value = 22939.0005
$ws.Range("M85").Value = -21535.0005
$ws.Range("H99").Value = 5998.1113
$ws.Range("I99").Value = 218.5
$ws.Range("J99").Value = 7649.4287
$ws.Range("K99").Value = 655.5
$ws.Range("L99").Value = 22948.2861
$ws.Range("M99").Value = 842.5
$ws.Range("N99").Value = -25944.2861
$ws.Range("H101").Value = 735.3889
$ws.Range("I101").Value = 995.75
$ws.Range("J101").Value = 527.1
$ws.Range("K101").Value = 2987.25
$ws.Range("L101").Value = 1581.3
$ws.Range("M101").Value = -1365.25
$ws.Range("N101").Value = -4825.3
$ws.Range("H107").Value = 1149.409
$ws.Range("I107").Value = 405.82352
$ws.Range("K107").Value = 405.82352
$ws.Range("M107").Value = 1514.17648
$ws.Range("H112").Value = 3892.647
$ws.Range("J112").Value = 4073.4375
$ws.Range("L112").Value = 12220.3125
$ws.Range("N112").Value = -14436.3125
$ws.Range("H138").Value = 3587.8645
$ws.Range("I138").Value = 2301.0454
$ws.Range("J138").Value = 4353
$ws.Range("K138").Value = 6903.1362
$ws.Range("L138").Value = 13059
$ws.Range("M138").Value = -1763.1362
$ws.Range("N138").Value = -23339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 100000
$ws.Range("J7").Value = 100000
$ws.Range("L7").Value = 100000
$ws.Range("N7").Value = -100228
$ws.Range("H37").Value = 750022500
$ws.Range("J37").Value = 750022500
$ws.Range("L37").Value = 750022500
$ws.Range("N37").Value = -750023046
$ws.Range("H53").Value = 250014940
$ws.Range("I53").Value = 9894.5
$ws.Range("J53").Value = 500020000
$ws.Range("K53").Value = 9894.5
$ws.Range("L53").Value = 500020000
$ws.Range("M53").Value = -9212.5
$ws.Range("N53").Value = -500021364
$ws.Range("H61").Value = 3754535.2
$ws.Range("I61").Value = 4549931
$ws.Range("J61").Value = 838084.5600000001
$ws.Range("K61").Value = 4549931
$ws.Range("L61").Value = 838084.5600000001
$ws.Range("M61").Value = -4549719
$ws.Range("N61").Value = -838508.5600000001
$ws.Range("H74").Value = 1373.1111
$ws.Range("I74").Value = 1359.9062
$ws.Range("J74").Value = 1478.75
$ws.Range("K74").Value = 1359.9062
$ws.Range("L74").Value = 1478.75
$ws.Range("M74").Value = -485.9061999999999
$ws.Range("N74").Value = -3226.75
$ws.Range("H77").Value = 1373.1111
$ws.Range("I77").Value = 1359.9062
$ws.Range("J77").Value = 1478.75
$ws.Range("K77").Value = 6799.530999999999
$ws.Range("L77").Value = 7393.75
$ws.Range("M77").Value = -2431.530999999999
$ws.Range("N77").Value = -16129.75
$ws.Range("H110").Value = 6706.75
$ws.Range("J110").Value = 5372.143
$ws.Range("L110").Value = 5372.143
$ws.Range("N110").Value = -9462.143
$ws.Range("H132").Value = 1962745
$ws.Range("I132").Value = 1880.4894
$ws.Range("K132").Value = 5641.468199999999
$ws.Range("M132").Value = -3111.468199999999
$ws.Range("H136").Value = 3754535.2
$ws.Range("I136").Value = 4549931
$ws.Range("J136").Value = 838084.5600000001
$ws.Range("K136").Value = 13649793
$ws.Range("L136").Value = 2514253.68
$ws.Range("M136").Value = -13647243
$ws.Range("N136").Value = -2519353.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2826.5405
$ws.Range("I20").Value = 2457.5833
$ws.Range("J20").Value = 3507.6924
$ws.Range("K20").Value = 2457.5833
$ws.Range("L20").Value = 3507.6924
$ws.Range("M20").Value = -2210.5833
$ws.Range("N20").Value = -4001.6924
$ws.Range("H43").Value = 295242
$ws.Range("J43").Value = 295242
$ws.Range("L43").Value = 295242
$ws.Range("N43").Value = -295604
$ws.Range("H99").Value = 1680.8462
$ws.Range("I99").Value = 650.1111
$ws.Range("K99").Value = 650.1111
$ws.Range("M99").Value = 847.8889
$ws.Range("H134").Value = 3849749.2
$ws.Range("I134").Value = 3598.5908
$ws.Range("K134").Value = 10795.7724
$ws.Range("M134").Value = -8260.7724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 12181.777
$ws.Range("I99").Value = 6274.3335
$ws.Range("J99").Value = 23996.666
$ws.Range("K99").Value = 6274.3335
$ws.Range("L99").Value = 23996.666
$ws.Range("M99").Value = -4776.3335
$ws.Range("N99").Value = -26992.666
$ws.Range("H126").Value = 12181.777
$ws.Range("I126").Value = 6274.3335
$ws.Range("J126").Value = 23996.666
$ws.Range("K126").Value = 18823.0005
$ws.Range("L126").Value = 71989.99800000001
$ws.Range("M126").Value = -16353.0005
$ws.Range("N126").Value = -76929.99800000001
$ws.Range("H141").Value = 547543.2
$ws.Range("I141").Value = 150000
$ws.Range("K141").Value = 150000
$ws.Range("M141").Value = -144820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 4620374.5
$ws.Range("J33").Value = 7700557.5
$ws.Range("L33").Value = 46203345
$ws.Range("N33").Value = -46203911
$ws.Range("H63").Value = 29168
$ws.Range("I63").Value = 20012
$ws.Range("K63").Value = 60036
$ws.Range("M63").Value = -59287
$ws.Range("H66").Value = 29168
$ws.Range("I66").Value = 20012
$ws.Range("K66").Value = 180108
$ws.Range("M66").Value = -176364
$ws.Range("H133").Value = 17107.092
$ws.Range("I133").Value = 10094.571
$ws.Range("K133").Value = 30283.713
$ws.Range("M133").Value = -25223.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2921.7646
$ws.Range("I102").Value = 2791.875
$ws.Range("K102").Value = 2791.875
$ws.Range("M102").Value = -1169.875
$ws.Range("H106").Value = 54500
$ws.Range("J106").Value = 54500
$ws.Range("L106").Value = 54500
$ws.Range("N106").Value = -57024
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H132").Value = 1670065.9
$ws.Range("I132").Value = 3377.327
$ws.Range("K132").Value = 10131.981
$ws.Range("M132").Value = -7601.981

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7774039.5
$ws.Range("I22").Value = 22020130
$ws.Range("J22").Value = 3444.2727
$ws.Range("K22").Value = 22020130
$ws.Range("L22").Value = 3444.2727
$ws.Range("M22").Value = -22019835
$ws.Range("N22").Value = -4034.2727
$ws.Range("H27").Value = 7774039.5
$ws.Range("I27").Value = 22020130
$ws.Range("J27").Value = 3444.2727
$ws.Range("K27").Value = 22020130
$ws.Range("L27").Value = 3444.2727
$ws.Range("M27").Value = -22020023
$ws.Range("N27").Value = -3658.2727
$ws.Range("H40").Value = 5955.3687
$ws.Range("J40").Value = 8658.333000000001
$ws.Range("L40").Value = 8658.333000000001
$ws.Range("N40").Value = -8930.333000000001
$ws.Range("H55").Value = 1294.2174
$ws.Range("I55").Value = 961.8333
$ws.Range("J55").Value = 1656.8182
$ws.Range("K55").Value = 961.8333
$ws.Range("L55").Value = 1656.8182
$ws.Range("M55").Value = -788.8333
$ws.Range("N55").Value = -2002.8182
$ws.Range("H68").Value = 3791546
$ws.Range("I68").Value = 6945836
$ws.Range("J68").Value = 6398.2
$ws.Range("K68").Value = 6945836
$ws.Range("L68").Value = 6398.2
$ws.Range("M68").Value = -6945087
$ws.Range("N68").Value = -7896.2
$ws.Range("H71").Value = 3791546
$ws.Range("I71").Value = 6945836
$ws.Range("J71").Value = 6398.2
$ws.Range("K71").Value = 34729180
$ws.Range("L71").Value = 31991
$ws.Range("M71").Value = -34725436
$ws.Range("N71").Value = -39479
$ws.Range("H82").Value = 2214
$ws.Range("I82").Value = 782.5714
$ws.Range("J82").Value = 3884
$ws.Range("K82").Value = 782.5714
$ws.Range("L82").Value = 3884
$ws.Range("M82").Value = -421.5714
$ws.Range("N82").Value = -4606
$ws.Range("H85").Value = 2214
$ws.Range("I85").Value = 782.5714
$ws.Range("J85").Value = 3884
$ws.Range("K85").Value = 782.5714
$ws.Range("L85").Value = 3884
$ws.Range("M85").Value = 465.4286
$ws.Range("N85").Value = -6380
$ws.Range("H93").Value = 2318099
$ws.Range("I93").Value = 2353.4211
$ws.Range("K93").Value = 2353.4211
$ws.Range("M93").Value = -1105.4211
$ws.Range("H136").Value = 3979.4075
$ws.Range("I136").Value = 2278.0625
$ws.Range("J136").Value = 6454.091
$ws.Range("K136").Value = 6834.1875
$ws.Range("L136").Value = 19362.273
$ws.Range("M136").Value = -4284.1875
$ws.Range("N136").Value = -24462.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6355.857
$ws.Range("I126").Value = 7181.8335
$ws.Range("J126").Value = 1400
$ws.Range("K126").Value = 21545.5005
$ws.Range("L126").Value = 4200
$ws.Range("M126").Value = -19075.5005
$ws.Range("N126").Value = -9140
$ws.Range("H132").Value = 530100.8
$ws.Range("I132").Value = 4037.2144
$ws.Range("K132").Value = 12111.6432
$ws.Range("M132").Value = -9581.643199999999
$ws.Range("H136").Value = 314456.38
$ws.Range("I136").Value = 2121.318
$ws.Range("K136").Value = 6363.954000000001
$ws.Range("M136").Value = -3813.954000000001
